$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.3082
$ws.Range("H2").Value = 6.9246
$ws.Range("I2").Value = 0.6638288620319053
$ws.Range("J2").Value = 0.6638288620319053
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 48.7014996656
$ws.Range("R2").Value = 438.3134969904
$ws.Range("S2").Value = 0.1936545589678796
$ws.Range("T2").Value = 0.1936545589678796
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.3082
$ws.Range("H3").Value = 6.9246
$ws.Range("I3").Value = 0.6638288620319053
$ws.Range("J3").Value = 0.6638288620319053
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 82.67842602219999
$ws.Range("R3").Value = 744.1058341998
$ws.Range("S3").Value = 0.3287589548047719
$ws.Range("T3").Value = 0.3287589548047719
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.3082
$ws.Range("H4").Value = 6.9246
$ws.Range("I4").Value = 0.6638288620319053
$ws.Range("J4").Value = 0.6638288620319053
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 35.56404544599999
$ws.Range("R4").Value = 320.076409014
$ws.Range("S4").Value = 0.1414153482592538
$ws.Range("T4").Value = 0.1414153482592538
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.168901
$ws.Range("H5").Value = 3.506703
$ws.Range("I5").Value = 0.3361711379680947
$ws.Range("J5").Value = 0.3361711379680947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 24.66304118387467
$ws.Range("R5").Value = 221.967370654872
$ws.Range("S5").Value = 0.09806906144706415
$ws.Range("T5").Value = 0.09806906144706416
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.168901
$ws.Range("H6").Value = 3.506703
$ws.Range("I6").Value = 0.3361711379680947
$ws.Range("J6").Value = 0.3361711379680947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 41.86937650800433
$ws.Range("R6").Value = 376.824388572039
$ws.Range("S6").Value = 0.1664875968418043
$ws.Range("T6").Value = 0.1664875968418043
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.168901
$ws.Range("H7").Value = 3.506703
$ws.Range("I7").Value = 0.3361711379680947
$ws.Range("J7").Value = 0.3361711379680947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 18.01007204136333
$ws.Range("R7").Value = 162.09064837227
$ws.Range("S7").Value = 0.07161447967922627
$ws.Range("T7").Value = 0.07161447967922627
